# SRS_Review.xlsx - "new update of SRS document after viewing the review sheet"
# Mark the open review points as Accepted (Decision column, E) on the
# "Cross review points " sheet, and highlight the last point (G10) with a
# bordered callout box, then leave the selection where the reviewer left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross review points ")

# --- Decision column: mark the previously-blank rows as "Accepted" ---
$ws.Range("E7").Value = "Accepted"
$ws.Range("E8").Value = "Accepted"
$ws.Range("E10").Value = "Accepted"
$ws.Range("E11").Value = "Accepted"
$ws.Range("E12").Value = "Accepted"
$ws.Range("E13").Value = "Accepted"
$ws.Range("E14").Value = "Accepted"
$ws.Range("E15").Value = "Accepted"
$ws.Range("E16").Value = "Accepted"

# --- Emphasize the last empty "Comment" cell (G10) with a bordered box ---
$g10 = $ws.Range("G10")
$g10.HorizontalAlignment = 1
$g10.VerticalAlignment = -4108
$g10.WrapText = $true
$g10.Font.Name = "Arial"
$g10.Font.Size = 12
$g10.Borders.LineStyle = 1
$g10.Borders.Weight = -4138
$g10.Borders.ColorIndex = 1

# --- Leave the view/selection where the reviewer ended up ---
[void]$ws.Activate()
[void]$ws.Range("E12").Select()
